$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.139.69"

# Row 3
$ws.Range("D3").Value = "'1.679.96"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'215.21"

# Row 6
$ws.Range("E6").Value = "  +0.48%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  +2.13%  "

# Row 9
$ws.Range("D9").Value = "'21.40"
$ws.Range("E9").Value = "  +5.30%  "

# Row 10
$ws.Range("D10").Value = "'0.0622"
$ws.Range("E10").Value = "  +0.55%  "

# Row 11
$ws.Range("E11").Value = "  +0.17%  "

# Row 12
$ws.Range("D12").Value = "'1.917.23"
$ws.Range("E12").Value = "  +0.39%  "

# Row 13
$ws.Range("D13").Value = "'1.671.72"
$ws.Range("E13").Value = "  -0.08%  "

# Row 14
$ws.Range("D14").Value = "'4.14"
$ws.Range("E14").Value = "  +1.53%  "

# Row 15
$ws.Range("E15").Value = "  +2.08%  "

# Row 16
$ws.Range("D16").Value = "'66.24"
$ws.Range("E16").Value = "  +0.94%  "

# Row 17
$ws.Range("D17").Value = "'27.133.83"
$ws.Range("E17").Value = "  +0.66%  "

# Row 18
$ws.Range("D18").Value = "'239.32"
$ws.Range("E18").Value = "  +1.63%  "

# Row 19
$ws.Range("D19").Value = "'8.06"
$ws.Range("E19").Value = "  -0.59%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0742"
$ws.Range("E20").Value = "  +1.24%  "

# Row 21
$ws.Range("E21").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").Value = "'9.48"
$ws.Range("E23").Value = "  +3.07%  "

# Row 24
$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -3.81%  "

# Row 25
$ws.Range("D25").Value = "'147.09"
$ws.Range("E25").Value = "  +1.17%  "

# Row 26
$ws.Range("D26").Value = "'7.25"
$ws.Range("E26").Value = "  +0.37%  "

# Row 27
$ws.Range("D27").Value = "'16.35"
$ws.Range("E27").Value = "  +2.18%  "

# Row 28
$ws.Range("E28").Value = "  +0.12%  "

# Row 29
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +0.42%  "

# Row 32
$ws.Range("D32").Value = "'1.564.99"
$ws.Range("E32").Value = "  +5.75%  "

# Row 33
$ws.Range("E33").Value = "  +1.68%  "

# Row 34
$ws.Range("D34").Value = "'3.21"
$ws.Range("E34").Value = "  +2.65%  "

# Row 35
$ws.Range("E35").Value = "  +0.21%  "

# Row 36
$ws.Range("E36").Value = "  +2.13%  "

# Row 37
$ws.Range("E37").Value = "  -1.12%  "

# Row 38
$ws.Range("E38").Value = "  +3.99%  "

# Row 39
$ws.Range("E39").Value = "  +2.11%  "

# Row 40
$ws.Range("E40").Value = "  +2.10%  "

# Row 41
$ws.Range("D41").Value = "'69.05"
$ws.Range("E41").Value = "  +3.06%  "

# Row 42
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$ws.Range("E43").Value = "  -5.17%  "

# Row 44
$ws.Range("E44").Value = "  -2.34%  "

# Row 45
$ws.Range("D45").Value = "'1.824.96"
$ws.Range("E45").Value = "  +0.62%  "

# Row 46
$ws.Range("D46").Value = "'0.785"
$ws.Range("E46").Value = "  +1.26%  "

# Row 47
$ws.Range("D47").Value = "'90.63"
$ws.Range("E47").Value = "  +0.25%  "

# Row 48
$ws.Range("E48").Value = "  +3.40%  "

# Row 49
$ws.Range("E49").Value = "  +1.56%  "

# Row 50 and 51: rows swap content (EnergySwap/Algorand swap places with updated values)
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.15"
$ws.Range("E50").Value = "  +6.21%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.104"
$ws.Range("E51").Value = "  +1.87%  "
